# Actualización automática 2025-11-17 15:30:09
# Inserts a new client "SANCHEZ BONILLA MARCO VINICIO" (alphabetically before
# "SANCHEZ SARMIENTO ANDRES FERNANDO") for asesor "ALMEIDA CUATIN JHONATHANN
# CARLOS" into the "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting the
# remaining client rows down by one. Also records a new PORCELANATO sale
# (59.02, in noviembre) for client "TIERRA GUAÑO JAIRO GABRIEL" and updates
# the dependent totals on "VENTAS POR GRUPO", "VENTA MENSUAL" and
# "CUMPLIMIENTO MENSUAL".

$wb = $excel.ActiveWorkbook

$asesor = "ALMEIDA CUATIN JHONATHANN CARLOS"

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO (product-by-client matrix)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert the new client row before row 29 (SANCHEZ SARMIENTO ...), shifting
# rows 29-37 down to 30-38.
$ws1.Rows("29:29").Insert()

$ws1.Range("A29").Value2 = $asesor
$ws1.Range("B29").Value2 = "SANCHEZ BONILLA MARCO VINICIO"
$ws1.Range("C29:R29").Value2 = 0

# New data point: TIERRA GUAÑO JAIRO GABRIEL (now row 34) has a PORCELANATO
# (column M) sale of 59.02.
$ws1.Range("M34").Value2 = 59.02

# Update the "X de 35" -> "X de 36" summary/count row, now on row 38.
$ws1.Range("C38").Value2 = "0 de 36"
$ws1.Range("D38").Value2 = "0 de 36"
$ws1.Range("E38").Value2 = "0 de 36"
$ws1.Range("F38").Value2 = "0 de 36"
$ws1.Range("G38").Value2 = "0 de 36"
$ws1.Range("H38").Value2 = "1 de 36"
$ws1.Range("I38").Value2 = "1 de 36"
$ws1.Range("J38").Value2 = "0 de 36"
$ws1.Range("K38").Value2 = "0 de 36"
$ws1.Range("L38").Value2 = "0 de 36"
$ws1.Range("M38").Value2 = "5 de 36"
$ws1.Range("N38").Value2 = "0 de 36"
$ws1.Range("O38").Value2 = "0 de 36"
$ws1.Range("P38").Value2 = "0 de 36"
$ws1.Range("Q38").Value2 = "0 de 36"
$ws1.Range("R38").Value2 = "0 de 36"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL (month-by-client matrix)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same client insertion as sheet 1.
$ws2.Rows("29:29").Insert()

$ws2.Range("A29").Value2 = $asesor
$ws2.Range("B29").Value2 = "SANCHEZ BONILLA MARCO VINICIO"
$ws2.Range("C29:G29").Value2 = 0

# TIERRA GUAÑO JAIRO GABRIEL (now row 34): noviembre (F) sale of 59.02, and
# PRESUPUESTO (G) of 1700.
$ws2.Range("F34").Value2 = 59.02
$ws2.Range("G34").Value2 = 1700

# Totals row (now row 38): noviembre total increases by the same 59.02.
$ws2.Range("F38").Value2 = 1345.18

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL (compliance summary)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# PORCELANATO row (12): VENTA (D) up by 59.02, POR CUMPLIR (E) down by the
# same amount, CUMPLIMIENTO (F) recalculated as VENTA / PRESUPUESTO.
$ws3.Range("D12").Value2 = 2015.2
$ws3.Range("E12").Value2 = 29198.8
$ws3.Range("F12").Value2 = 0.06456077401166144

# TOTAL row (14): same 59.02 ripple.
$ws3.Range("D14").Value2 = 2104.84
$ws3.Range("E14").Value2 = 38174.72164865473
$ws3.Range("F14").Value2 = 0.05225578218451884
